$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.544.75"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.328.82"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "2.679.91"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "2.325.52"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "43.205.87"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0889"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("E34").Value = "  +9.34%  "
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0364"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  +5.77%  "
$ws.Range("E41").Value = "  +8.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.18%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "1.666.95"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.74%  "
$ws.Range("E51").Value = "  +10.17%  "
